$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp banner (A1)
$ws.Range("A1").Value = "Datos actualizados a 22 de Octubre de 2020 a las 21:58"

# Row 4: Estados Unidos
$ws.Cells.Item(4, 1).Value = "Estados Unidos"
$ws.Cells.Item(4, 2).Value = 8635113
$ws.Cells.Item(4, 3).Value = 50294
$ws.Cells.Item(4, 4).Value = 5631598
$ws.Cells.Item(4, 5).Value = 2775498
$ws.Cells.Item(4, 6).Value = 0
$ws.Cells.Item(4, 7).Value = 609
$ws.Cells.Item(4, 8).Value = 228017

# Row 5: India
$ws.Cells.Item(5, 1).Value = "India"
$ws.Cells.Item(5, 2).Value = 7759252
$ws.Cells.Item(5, 3).Value = 54094
$ws.Cells.Item(5, 4).Value = 6945879
$ws.Cells.Item(5, 5).Value = 696043
$ws.Cells.Item(5, 6).Value = 0
$ws.Cells.Item(5, 7).Value = 677
$ws.Cells.Item(5, 8).Value = 117330

# Row 15: Sudafrica
$ws.Cells.Item(15, 1).Value = "Sudafrica"
$ws.Cells.Item(15, 2).Value = 710515
$ws.Cells.Item(15, 3).Value = 2156
$ws.Cells.Item(15, 4).Value = 642560
$ws.Cells.Item(15, 5).Value = 49112
$ws.Cells.Item(15, 6).Value = 0
$ws.Cells.Item(15, 7).Value = 102
$ws.Cells.Item(15, 8).Value = 18843

# Row 20: Alemania
$ws.Cells.Item(20, 1).Value = "Alemania"
$ws.Cells.Item(20, 2).Value = 403844
$ws.Cells.Item(20, 3).Value = 12489
$ws.Cells.Item(20, 4).Value = 302100
$ws.Cells.Item(20, 5).Value = 91700
$ws.Cells.Item(20, 6).Value = 0
$ws.Cells.Item(20, 7).Value = 45
$ws.Cells.Item(20, 8).Value = 10044

# Row 28: Israel
$ws.Cells.Item(28, 1).Value = "Israel"
$ws.Cells.Item(28, 2).Value = 308247
$ws.Cells.Item(28, 3).Value = 912
$ws.Cells.Item(28, 4).Value = 288337
$ws.Cells.Item(28, 5).Value = 17591
$ws.Cells.Item(28, 6).Value = 0
$ws.Cells.Item(28, 7).Value = 28
$ws.Cells.Item(28, 8).Value = 2319

# Row 53: Etiopia
$ws.Cells.Item(53, 1).Value = "Etiopia"
$ws.Cells.Item(53, 2).Value = 91693
$ws.Cells.Item(53, 3).Value = 575
$ws.Cells.Item(53, 4).Value = 45260
$ws.Cells.Item(53, 5).Value = 45037
$ws.Cells.Item(53, 6).Value = 0
$ws.Cells.Item(53, 7).Value = 12
$ws.Cells.Item(53, 8).Value = 1396

# Row 100: Zambia
$ws.Cells.Item(100, 1).Value = "Zambia"
$ws.Cells.Item(100, 2).Value = 16035
$ws.Cells.Item(100, 3).Value = 35
$ws.Cells.Item(100, 4).Value = 15168
$ws.Cells.Item(100, 5).Value = 521
$ws.Cells.Item(100, 6).Value = 0
$ws.Cells.Item(100, 7).Value = 0
$ws.Cells.Item(100, 8).Value = 346

# Row 117: Angola
$ws.Cells.Item(117, 1).Value = "Angola"
$ws.Cells.Item(117, 2).Value = 8582
$ws.Cells.Item(117, 3).Value = 244
$ws.Cells.Item(117, 4).Value = 3305
$ws.Cells.Item(117, 5).Value = 5017
$ws.Cells.Item(117, 6).Value = 0
$ws.Cells.Item(117, 7).Value = 5
$ws.Cells.Item(117, 8).Value = 260

# Row 138: Ruanda
$ws.Cells.Item(138, 1).Value = "Ruanda"
$ws.Cells.Item(138, 2).Value = 5017
$ws.Cells.Item(138, 3).Value = 5
$ws.Cells.Item(138, 4).Value = 4803
$ws.Cells.Item(138, 5).Value = 180
$ws.Cells.Item(138, 6).Value = 0
$ws.Cells.Item(138, 7).Value = 0
$ws.Cells.Item(138, 8).Value = 34

# Row 139: Reunion
$ws.Cells.Item(139, 1).Value = "Reunion"
$ws.Cells.Item(139, 2).Value = 5015
$ws.Cells.Item(139, 3).Value = 0
$ws.Cells.Item(139, 4).Value = 4445
$ws.Cells.Item(139, 5).Value = 551
$ws.Cells.Item(139, 6).Value = 0
$ws.Cells.Item(139, 7).Value = 0
$ws.Cells.Item(139, 8).Value = 19

# Row 141: Aruba
$ws.Cells.Item(141, 1).Value = "Aruba"
$ws.Cells.Item(141, 2).Value = 4389
$ws.Cells.Item(141, 3).Value = 20
$ws.Cells.Item(141, 4).Value = 4120
$ws.Cells.Item(141, 5).Value = 233
$ws.Cells.Item(141, 6).Value = 0
$ws.Cells.Item(141, 7).Value = 1
$ws.Cells.Item(141, 8).Value = 36

# Row 180: Comoras
$ws.Cells.Item(180, 1).Value = "Comoras"
$ws.Cells.Item(180, 2).Value = 517
$ws.Cells.Item(180, 3).Value = 13
$ws.Cells.Item(180, 4).Value = 494
$ws.Cells.Item(180, 5).Value = 16
$ws.Cells.Item(180, 6).Value = 0
$ws.Cells.Item(180, 7).Value = 0
$ws.Cells.Item(180, 8).Value = 7

# Row 181: Tanzania
$ws.Cells.Item(181, 1).Value = "Tanzania"
$ws.Cells.Item(181, 2).Value = 509
$ws.Cells.Item(181, 3).Value = 0
$ws.Cells.Item(181, 4).Value = 183
$ws.Cells.Item(181, 5).Value = 305
$ws.Cells.Item(181, 6).Value = 0
$ws.Cells.Item(181, 7).Value = 0
$ws.Cells.Item(181, 8).Value = 21

# Row 192: Barbados
$ws.Cells.Item(192, 1).Value = "Barbados"
$ws.Cells.Item(192, 2).Value = 224
$ws.Cells.Item(192, 3).Value = 2
$ws.Cells.Item(192, 4).Value = 207
$ws.Cells.Item(192, 5).Value = 10
$ws.Cells.Item(192, 6).Value = 0
$ws.Cells.Item(192, 7).Value = 0
$ws.Cells.Item(192, 8).Value = 7

# Row 216: Islas Malvinas
$ws.Cells.Item(216, 1).Value = "Islas Malvinas"
$ws.Cells.Item(216, 2).Value = 13
$ws.Cells.Item(216, 3).Value = 0
$ws.Cells.Item(216, 4).Value = 13
$ws.Cells.Item(216, 5).Value = 0
$ws.Cells.Item(216, 6).Value = 0
$ws.Cells.Item(216, 7).Value = 0
$ws.Cells.Item(216, 8).Value = 0

# Row 217: Montserrat
$ws.Cells.Item(217, 1).Value = "Montserrat"
$ws.Cells.Item(217, 2).Value = 13
$ws.Cells.Item(217, 3).Value = 0
$ws.Cells.Item(217, 4).Value = 12
$ws.Cells.Item(217, 5).Value = 0
$ws.Cells.Item(217, 6).Value = 0
$ws.Cells.Item(217, 7).Value = 0
$ws.Cells.Item(217, 8).Value = 1
